$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text storage for Price (D) cells so numeric-looking strings
# (e.g. "3.00", "14.95") keep their exact text representation instead
# of being parsed/rounded as numbers.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

# Price (column D) updates
$ws.Range("D2").Value = "41.249.62"
$ws.Range("D3").Value = "2.462.25"
$ws.Range("D4").Value = "0.999"
$ws.Range("D5").Value = "311.97"
$ws.Range("D6").Value = "94.09"
$ws.Range("D7").Value = "0.547"
$ws.Range("D14").Value = "2.841.37"
$ws.Range("D15").Value = "2.460.00"
$ws.Range("D16").Value = "14.95"
$ws.Range("D17").Value = "0.785"
$ws.Range("D18").Value = "41.239.81"
$ws.Range("D20").Value = "0.0₃0922"
$ws.Range("D21").Value = "11.19"
$ws.Range("D22").Value = "68.41"
$ws.Range("D23").Value = "236.69"
$ws.Range("D24").Value = "2.74"
$ws.Range("D26").Value = "1.89"
$ws.Range("D27").Value = "23.94"
$ws.Range("D30").Value = "36.43"
$ws.Range("D31").Value = "151.43"
$ws.Range("D32").Value = "5.47"
$ws.Range("D33").Value = "2.61"
$ws.Range("D35").Value = "0.0743"
$ws.Range("D36").Value = "3.03"
$ws.Range("D37").Value = "17.11"
$ws.Range("D38").Value = "1.86"
$ws.Range("D41").Value = "4.24"
$ws.Range("D43").Value = "19.49"
$ws.Range("D44").Value = "1.980.23"
$ws.Range("D45").Value = "0.0285"
$ws.Range("D46").Value = "3.00"
$ws.Range("D47").Value = "8.69"
$ws.Range("D48").Value = "2.704.65"
$ws.Range("D49").Value = "69.36"
$ws.Range("D50").Value = "96.27"
$ws.Range("D51").Value = "74.39"

# Volume(1h) (column E) updates
$ws.Range("E2").Value = "  -3.77%  "
$ws.Range("E3").Value = "  -3.01%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("E5").Value = "  +0.25%  "
$ws.Range("E6").Value = "  -6.72%  "
$ws.Range("E7").Value = "  -3.64%  "
$ws.Range("E8").Value = "  +0.05%  "
$ws.Range("E9").Value = "  -5.14%  "
$ws.Range("E10").Value = "  -6.86%  "
$ws.Range("E11").Value = "  -3.29%  "
$ws.Range("E12").Value = "  -0.92%  "
$ws.Range("E13").Value = "  -4.94%  "
$ws.Range("E14").Value = "  -2.94%  "
$ws.Range("E15").Value = "  -3.25%  "
$ws.Range("E16").Value = "  -2.89%  "
$ws.Range("E17").Value = "  -4.09%  "
$ws.Range("E19").Value = "  -6.11%  "
$ws.Range("E20").Value = "  -3.48%  "
$ws.Range("E21").Value = "  -9.76%  "
$ws.Range("E22").Value = "  -2.53%  "
$ws.Range("E23").Value = "  -3.12%  "
$ws.Range("E24").Value = "  -5.15%  "
$ws.Range("E25").Value = "  +0.10%  "
$ws.Range("E26").Value = "  -7.15%  "
$ws.Range("E27").Value = "  -6.31%  "
$ws.Range("E28").Value = "  -4.72%  "
$ws.Range("E29").Value = "  -5.68%  "
$ws.Range("E30").Value = "  -6.01%  "
$ws.Range("E31").Value = "  -4.66%  "
$ws.Range("E32").Value = "  -7.02%  "
$ws.Range("E33").Value = "  -6.03%  "
$ws.Range("E34").Value = "  -3.14%  "
$ws.Range("E35").Value = "  -6.44%  "
$ws.Range("E36").Value = "  -3.50%  "
$ws.Range("E37").Value = "  -6.36%  "
$ws.Range("E38").Value = "  -5.75%  "
$ws.Range("E39").Value = "  -3.30%  "
$ws.Range("E40").Value = "  -8.73%  "
$ws.Range("E41").Value = "  +1.61%  "
$ws.Range("E42").Value = "  +0.17%  "
$ws.Range("E43").Value = "  -11.28%  "
$ws.Range("E44").Value = "  -0.95%  "
$ws.Range("E45").Value = "  -5.12%  "
$ws.Range("E46").Value = "  -9.34%  "
$ws.Range("E47").Value = "  -4.05%  "
$ws.Range("E48").Value = "  -2.70%  "
$ws.Range("E49").Value = "  -4.56%  "
$ws.Range("E50").Value = "  -5.18%  "
$ws.Range("E51").Value = "  -7.48%  "
